# Negative_Manifest.xlsx update — add rows n1..n3 (subjects with no meltpatch,
# classified negative), matching the 2021-02-06 manifest commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("n1", "n1_IMG_3178.jpeg", "True", "no_meltpatch", "negative"),
    @("n2", "n2_IMG_3180.jpeg", "True", "no_meltpatch", "negative"),
    @("n3", "n3_IMG_3174.jpeg", "True", "no_meltpatch", "negative")
)

$row = 2
foreach ($rec in $data) {
    $col = 1
    foreach ($val in $rec) {
        $cell = $ws.Cells.Item($row, $col)
        if ($val -eq "True" -or $val -eq "False") {
            # Force literal text (not a Boolean) the same way typing `'True`
            # into the cell in the Excel UI does, then drop the resulting
            # quote-prefix style so no extra formatting sticks around.
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
        $col++
    }
    $row++
}

# Column B ("#file_name") needs to widen to fit the new, longer file names.
$ws.Columns.Item(2).ColumnWidth = 22.36328125
